# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the per-language handback sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 06:47:45"
$wsZhCn.Range("H2").Value = "2016-03-14 06:48:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 06:47:48"
$wsDeDe.Range("H2").Value = "2016-03-14 06:48:09"
